$excel = New-Object -ComObject Excel.Application
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new record as row 726, pushing all existing rows
# (726..791) down by one (to 727..792). This matches the diff, which
# shows every row from 726 onward taking on the values that used to
# belong to the row just above it, with one new row of real data
# introduced at 726 and the dimension growing to A1:R792.
$ws.Rows.Item(726).Insert()

$ws.Range("A726").Value = 6
$ws.Range("B726").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C726").Value = "Metropolitana"
$ws.Range("D726").Value = 45132
$ws.Range("E726").Value = 13
$ws.Range("F726").Value = 100112012
$ws.Range("G726").Value = "Espinaca"
$ws.Range("H726").Value = "Sin especificar"
$ws.Range("I726").Value = "Primera"
$ws.Range("J726").Value = 520
$ws.Range("K726").Value = 6000
$ws.Range("L726").Value = 7000
$ws.Range("M726").Value = 6423
$ws.Range("N726").Value = "$/cuna 10 kilos"
$ws.Range("O726").Value = "Región Metropolitana"
$ws.Range("P726").Value = 642
$ws.Range("Q726").Value = 10
$ws.Range("R726").Value = "Hortaliza"
